# patient_log.xlsx -- add a 4th stent time/location pair and backfill the
# UI7616 procedure row (row 4) with the full set of LHC/stent/balloon data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new (blank) columns before Time_Balloon_1 (old col L) so
#     there's room for Time_Stent_4 / Location_Stent_4 after Location_Stent_3
#     (col K). Headers/labels for them are filled in further below, after
#     the row 4 data, to mirror the original authoring order. ---
$ws.Columns("L:M").Insert()

# --- Backfill row 4 (patient UI7616) with the LHC / stent / balloon data ---
# Time_LHC
$ws.Range("E4").Value2 = 0.65833333333333333
$ws.Range("E4").NumberFormat = "h:mm"

# Stent 1
$ws.Range("F4").Value2 = 0.71597222222222223
$ws.Range("F4").NumberFormat = "h:mm"
$ws.Range("G4").Value2 = "Marginal"

# Stent 2
$ws.Range("H4").Value2 = 0.74236111111111114
$ws.Range("H4").NumberFormat = "h:mm"
$ws.Range("I4").Value2 = "LCX"

# Stent 3
$ws.Range("J4").Value2 = 0.74513888888888891
$ws.Range("J4").NumberFormat = "h:mm"
$ws.Range("K4").Value2 = "LCX"

# New header labels for the inserted Stent 4 columns.
$ws.Range("L1").Value2 = "Time_Stent_4"
$ws.Range("M1").Value2 = "Location_Stent_4"

# Stent 4 (new)
$ws.Range("L4").Value2 = 0.75624999999999998
$ws.Range("L4").NumberFormat = "h:mm"
$ws.Range("M4").Value2 = "LCX"

# Balloon 1
$ws.Range("N4").Value2 = 0.68472222222222223
$ws.Range("N4").NumberFormat = "h:mm"
$ws.Range("O4").Value2 = "Unknown"

# Balloon 2
$ws.Range("P4").Value2 = 0.69305555555555554
$ws.Range("P4").NumberFormat = "h:mm"
$ws.Range("Q4").Value2 = "Marginal"

# Balloon 3
$ws.Range("R4").Value2 = 0.70694444444444438
$ws.Range("R4").NumberFormat = "h:mm"
$ws.Range("S4").Value2 = "LCX"

# Row 4 now wraps across three lines instead of two -- grow it to match.
$ws.Rows("4").RowHeight = 47.25

# Reflect the new working area: zoomed back to 100% and parked on the new
# Stent 4 columns.
$excel.ActiveWindow.Zoom = 100
$ws.Range("O5").Select() | Out-Null
